# Daily attendance processing - 2026-01-03 09:05:02
# Normalizes the "Recorded By" column (G): whenever the comma-separated
# list of recorders includes a "System" token, flip the order of the
# whole list (so "System" moves from the end to the front).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Trim() -eq "system") {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $count = $parts.Count
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newText = $reversed -join ", "
        $cell.Value = $newText
    }
}
